$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. New question in the "Probeneingang" bullet list, right after the
#    "Messung DSC, Pulver, IR ..." item: "Eingangsanalysen auflisten?"
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Messung DSC, Pulver, IR*") {
        $target = $p
    }
}

if ($target -ne $null) {
    # Inserting a paragraph break after the item duplicates its paragraph
    # formatting (style + numbering), exactly like pressing Enter at the end
    # of the line in Word.
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $newPara.Range.Text = "Eingangsanalysen auflisten?"
}

# ---------------------------------------------------------------------------
# 2. New character styles ListLabel 19 .. ListLabel 36 (list-bullet label
#    styles, matching the ListLabel 1..18 styles already in the template).
# ---------------------------------------------------------------------------
for ($i = 19; $i -le 36; $i++) {
    $styleId = "ListLabel$i"
    $style = $d.Styles.Add($styleId, 2)
    $style.NameLocal = "ListLabel $i"
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
}
